# Apply Betfair Back/Lay odds updates for 2025-11-17 workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 8).Value = 2.68  # H2: 2.52 -> 2.68
$ws.Cells.Item(2, 11).Value = 3.5  # K2: 3.9 -> 3.5
$ws.Cells.Item(3, 11).Value = 3.6  # K3: 3.7 -> 3.6
$ws.Cells.Item(3, 12).Value = 1.45  # L3: 1.01 -> 1.45
$ws.Cells.Item(4, 6).Value = 6  # F4: 5.9 -> 6
$ws.Cells.Item(4, 7).Value = 7  # G4: 6 -> 7
$ws.Cells.Item(4, 8).Value = 1.78  # H4: 1.74 -> 1.78
$ws.Cells.Item(4, 9).Value = 1.91  # I4: 1.89 -> 1.91
$ws.Cells.Item(4, 10).Value = 3.15  # J4: 3.2 -> 3.15
$ws.Cells.Item(4, 11).Value = 3.45  # K4: 3.6 -> 3.45
$ws.Cells.Item(4, 12).Value = 1.51  # L4: 1.6 -> 1.51
$ws.Cells.Item(4, 14).Value = 2.36  # N4: 2.38 -> 2.36
$ws.Cells.Item(4, 15).Value = 1.59  # O4: 1.6 -> 1.59
$ws.Cells.Item(4, 17).Value = 2.78  # Q4: 2.76 -> 2.78
$ws.Cells.Item(4, 20).Value = 2.42  # T4: 2.44 -> 2.42
$ws.Cells.Item(4, 21).Value = 1.56  # U4: 1.57 -> 1.56
$ws.Cells.Item(4, 22).Value = 2.1  # V4: 2.12 -> 2.1
$ws.Cells.Item(4, 23).Value = 1.17  # W4: 1.2 -> 1.17
$ws.Cells.Item(4, 24).Value = 8  # X4: 8.199999999999999 -> 8
$ws.Cells.Item(4, 26).Value = 9.199999999999999  # Z4: 9 -> 9.199999999999999
$ws.Cells.Item(4, 27).Value = 23  # AA4: 22 -> 23
$ws.Cells.Item(4, 28).Value = 14.5  # AB4: 15 -> 14.5
$ws.Cells.Item(4, 29).Value = 8.4  # AC4: 8.6 -> 8.4
$ws.Cells.Item(4, 32).Value = 50  # AF4: 55 -> 50
$ws.Cells.Item(4, 33).Value = 32  # AG4: 30 -> 32
$ws.Cells.Item(4, 41).Value = 26  # AO4: 24 -> 26
$ws.Cells.Item(5, 12).Value = 1.01  # L5: 1.45 -> 1.01
$ws.Cells.Item(6, 6).Value = 2.62  # F6: 2.36 -> 2.62
$ws.Cells.Item(6, 7).Value = 3  # G6: 3.3 -> 3
$ws.Cells.Item(6, 8).Value = 3.05  # H6: 2.66 -> 3.05
$ws.Cells.Item(6, 9).Value = 3.55  # I6: 3.85 -> 3.55
$ws.Cells.Item(6, 10).Value = 2.8  # J6: 2.64 -> 2.8
$ws.Cells.Item(6, 11).Value = 3.3  # K6: 3.7 -> 3.3
$ws.Cells.Item(6, 13).Value = 1.12  # M6: 1.01 -> 1.12
$ws.Cells.Item(6, 14).Value = 2.38  # N6: 1.38 -> 2.38
$ws.Cells.Item(6, 15).Value = 1.59  # O6: 1.01 -> 1.59
$ws.Cells.Item(6, 16).Value = 1.45  # P6: 1.38 -> 1.45
$ws.Cells.Item(6, 17).Value = 2.78  # Q6: 2.52 -> 2.78
$ws.Cells.Item(6, 19).Value = 5.2  # S6: 2.54 -> 5.2
$ws.Cells.Item(6, 20).Value = 2.16  # T6: 1.01 -> 2.16
$ws.Cells.Item(6, 21).Value = 1.71  # U6: 1.01 -> 1.71
$ws.Cells.Item(6, 22).Value = 1.39  # V6: 1.35 -> 1.39
$ws.Cells.Item(6, 23).Value = 1.5  # W6: 1.44 -> 1.5
$ws.Cells.Item(6, 24).Value = 9.199999999999999  # X6: 1000 -> 9.199999999999999
$ws.Cells.Item(6, 25).Value = 10.5  # Y6: 1000 -> 10.5
$ws.Cells.Item(6, 28).Value = 9.199999999999999  # AB6: 1000 -> 9.199999999999999
$ws.Cells.Item(6, 29).Value = 8.4  # AC6: 1000 -> 8.4
$ws.Cells.Item(6, 33).Value = 17  # AG6: 1000 -> 17
$ws.Cells.Item(7, 6).Value = 3.2  # F7: 3.3 -> 3.2
$ws.Cells.Item(7, 7).Value = 3.4  # G7: 3.5 -> 3.4
$ws.Cells.Item(7, 8).Value = 2.94  # H7: 2.78 -> 2.94
$ws.Cells.Item(7, 9).Value = 3.05  # I7: 2.98 -> 3.05
$ws.Cells.Item(7, 11).Value = 2.8  # K7: 2.84 -> 2.8
$ws.Cells.Item(7, 15).Value = 1.81  # O7: 1.83 -> 1.81
$ws.Cells.Item(7, 18).Value = 1.12  # R7: 1.11 -> 1.12
$ws.Cells.Item(7, 21).Value = 1.53  # U7: 1.52 -> 1.53
$ws.Cells.Item(7, 22).Value = 1.48  # V7: 1.51 -> 1.48
$ws.Cells.Item(7, 25).Value = 7  # Y7: 7.2 -> 7
$ws.Cells.Item(7, 27).Value = 160  # AA7: 150 -> 160
$ws.Cells.Item(7, 32).Value = 21  # AF7: 19.5 -> 21
$ws.Cells.Item(7, 33).Value = 22  # AG7: 18.5 -> 22
$ws.Cells.Item(7, 40).Value = 130  # AN7: 1000 -> 130
$ws.Cells.Item(8, 6).Value = 2.18  # F8: 2.2 -> 2.18
$ws.Cells.Item(8, 7).Value = 2.34  # G8: 2.36 -> 2.34
$ws.Cells.Item(8, 8).Value = 3.95  # H8: 3.8 -> 3.95
$ws.Cells.Item(8, 9).Value = 4.4  # I8: 4.2 -> 4.4
$ws.Cells.Item(8, 11).Value = 3.2  # K8: 3.3 -> 3.2
$ws.Cells.Item(8, 22).Value = 1.3  # V8: 1.31 -> 1.3
$ws.Cells.Item(8, 23).Value = 1.74  # W8: 1.73 -> 1.74
$ws.Cells.Item(8, 28).Value = 7.2  # AB8: 980 -> 7.2
$ws.Cells.Item(8, 29).Value = 7.6  # AC8: 980 -> 7.6
$ws.Cells.Item(9, 6).Value = 2.04  # F9: 2 -> 2.04
$ws.Cells.Item(9, 7).Value = 2.14  # G9: 2.12 -> 2.14
$ws.Cells.Item(9, 11).Value = 3.35  # K9: 3.4 -> 3.35
$ws.Cells.Item(9, 12).Value = 1.45  # L9: 1.46 -> 1.45
$ws.Cells.Item(9, 19).Value = 5.3  # S9: 4.4 -> 5.3
$ws.Cells.Item(9, 23).Value = 1.87  # W9: 1.89 -> 1.87
$ws.Cells.Item(9, 26).Value = 980  # Z9: 1000 -> 980
$ws.Cells.Item(9, 40).Value = 1000  # AN9: 32 -> 1000
$ws.Cells.Item(10, 7).Value = 2.46  # G10: 2.5 -> 2.46
$ws.Cells.Item(10, 9).Value = 3.55  # I10: 3.6 -> 3.55
$ws.Cells.Item(10, 10).Value = 3.2  # J10: 3.15 -> 3.2
$ws.Cells.Item(10, 11).Value = 3.35  # K10: 3.4 -> 3.35
$ws.Cells.Item(10, 14).Value = 2.54  # N10: 2.5 -> 2.54
$ws.Cells.Item(10, 15).Value = 1.62  # O10: 1.6 -> 1.62
$ws.Cells.Item(10, 16).Value = 1.51  # P10: 1.5 -> 1.51
$ws.Cells.Item(10, 18).Value = 1.17  # R10: 1.16 -> 1.17
$ws.Cells.Item(10, 20).Value = 2.26  # T10: 2.28 -> 2.26
$ws.Cells.Item(10, 21).Value = 1.71  # U10: 1.7 -> 1.71
$ws.Cells.Item(10, 23).Value = 1.68  # W10: 1.66 -> 1.68
$ws.Cells.Item(10, 24).Value = 8.199999999999999  # X10: 8.6 -> 8.199999999999999
$ws.Cells.Item(10, 25).Value = 9.6  # Y10: 8.800000000000001 -> 9.6
$ws.Cells.Item(10, 36).Value = 980  # AJ10: 38 -> 980
$ws.Cells.Item(10, 37).Value = 980  # AK10: 75 -> 980
$ws.Cells.Item(10, 41).Value = 90  # AO10: 95 -> 90
